# Update Thbs1-Tnfrsf11b.xlsx with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 41.428665
$ws.Range("H2").Value = 124.285995
$ws.Range("I2").Value = 0.06969137269740189
$ws.Range("J2").Value = 0.06969137269740189
$ws.Range("M2").Value = 0.08241233333333334
$ws.Range("Q2").Value = 3.414232949535
$ws.Range("R2").Value = 30.728096545815
$ws.Range("S2").Value = 0.003294315099086412
$ws.Range("T2").Value = 0.003294315099086412

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 41.428665
$ws.Range("H3").Value = 124.285995
$ws.Range("I3").Value = 0.06969137269740189
$ws.Range("J3").Value = 0.06969137269740189
$ws.Range("Q3").Value = 68.81400685295999
$ws.Range("R3").Value = 619.3260616766399
$ws.Range("S3").Value = 0.06639705759831548
$ws.Range("T3").Value = 0.06639705759831548

# Row 4 (FAPs -> ECs)
$ws.Range("H4").Value = 510.696747
$ws.Range("I4").Value = 0.2863649869040173
$ws.Range("J4").Value = 0.2863649869040173
$ws.Range("M4").Value = 0.08241233333333334
$ws.Range("S4").Value = 0.01353648900422299
$ws.Range("T4").Value = 0.01353648900422299

# Row 5 (FAPs -> FAPs)
$ws.Range("H5").Value = 510.696747
$ws.Range("I5").Value = 0.2863649869040173
$ws.Range("J5").Value = 0.2863649869040173
$ws.Range("Q5").Value = 282.759851162976
$ws.Range("S5").Value = 0.2728284978997943
$ws.Range("T5").Value = 0.2728284978997944

# Row 6 (Inflammatory-Mac -> ECs)
$ws.Range("G6").Value = 244.5761666666666
$ws.Range("H6").Value = 733.7284999999999
$ws.Range("I6").Value = 0.4114264551867299
$ws.Range("J6").Value = 0.41142645518673
$ws.Range("M6").Value = 0.08241233333333334
$ws.Range("Q6").Value = 20.15609257272222
$ws.Range("R6").Value = 181.4048331545
$ws.Range("S6").Value = 0.01944815162947381
$ws.Range("T6").Value = 0.01944815162947382

# Row 7 (Inflammatory-Mac -> FAPs)
$ws.Range("G7").Value = 244.5761666666666
$ws.Range("H7").Value = 733.7284999999999
$ws.Range("I7").Value = 0.4114264551867299
$ws.Range("J7").Value = 0.41142645518673
$ws.Range("Q7").Value = 406.2468826613333
$ws.Range("R7").Value = 3656.221943951999
$ws.Range("S7").Value = 0.3919783035572562
$ws.Range("T7").Value = 0.3919783035572562

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 24.173247
$ws.Range("H8").Value = 72.51974100000001
$ws.Range("I8").Value = 0.04066427836821081
$ws.Range("J8").Value = 0.04066427836821081
$ws.Range("M8").Value = 0.08241233333333334
$ws.Range("Q8").Value = 1.992173689513
$ws.Range("R8").Value = 17.929563205617
$ws.Range("S8").Value = 0.001922202720894949
$ws.Range("T8").Value = 0.001922202720894949

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 24.173247
$ws.Range("H9").Value = 72.51974100000001
$ws.Range("I9").Value = 0.04066427836821081
$ws.Range("J9").Value = 0.04066427836821081
$ws.Range("Q9").Value = 40.15234342492801
$ws.Range("R9").Value = 361.3710908243521
$ws.Range("S9").Value = 0.03874207564731587
$ws.Range("T9").Value = 0.03874207564731586

# Row 10 (Resolving-Mac -> ECs)
$ws.Range("G10").Value = 114.0486906666667
$ws.Range("H10").Value = 342.146072
$ws.Range("I10").Value = 0.19185290684364
$ws.Range("J10").Value = 0.19185290684364
$ws.Range("M10").Value = 0.08241233333333334
$ws.Range("Q10").Value = 9.399018711451557
$ws.Range("R10").Value = 84.591168403064
$ws.Range("S10").Value = 0.0090688976749368
$ws.Range("T10").Value = 0.0090688976749368

# Row 11 (Resolving-Mac -> FAPs)
$ws.Range("G11").Value = 114.0486906666667
$ws.Range("H11").Value = 342.146072
$ws.Range("I11").Value = 0.19185290684364
$ws.Range("J11").Value = 0.19185290684364
$ws.Range("Q11").Value = 189.4376123659094
$ws.Range("R11").Value = 1704.938511293184
$ws.Range("S11").Value = 0.1827840091687032
$ws.Range("T11").Value = 0.1827840091687032
